$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H76").Value = 3324.5881
$ws.Range("I76").Value = 2914.1025
$ws.Range("J76").Value = 4658.6665
$ws.Range("K76").Value = 2914.1025
$ws.Range("L76").Value = 4658.6665
$ws.Range("M76").Value = -2599.1025
$ws.Range("N76").Value = -5288.6665
$ws.Range("H79").Value = 3324.5881
$ws.Range("I79").Value = 2914.1025
$ws.Range("J79").Value = 4658.6665
$ws.Range("K79").Value = 2914.1025
$ws.Range("L79").Value = 4658.6665
$ws.Range("M79").Value = -1822.1025
$ws.Range("N79").Value = -6842.6665
$ws.Range("H114").Value = 39826.668
$ws.Range("J114").Value = 39826.668
$ws.Range("L114").Value = 39826.668
$ws.Range("N114").Value = -48504.668
$ws.Range("H129").Value = 1044.5209
$ws.Range("J129").Value = 1127.4524
$ws.Range("L129").Value = 3382.357199999999
$ws.Range("N129").Value = -13382.3572

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H13").Value = 5250000
$ws.Range("I13").Value = 5250000
$ws.Range("K13").Value = 5250000
$ws.Range("M13").Value = -5249856
$ws.Range("H61").Value = 1456.1666
$ws.Range("I61").Value = 1411.5
$ws.Range("K61").Value = 1411.5
$ws.Range("M61").Value = -1199.5
$ws.Range("H102").Value = 3257.85
$ws.Range("I102").Value = 3042.2666
$ws.Range("K102").Value = 3042.2666
$ws.Range("M102").Value = -1420.2666
$ws.Range("H132").Value = 1872.8704
$ws.Range("I132").Value = 965.96295
$ws.Range("J132").Value = 2779.7778
$ws.Range("K132").Value = 2897.88885
$ws.Range("L132").Value = 8339.3334
$ws.Range("M132").Value = -367.8888499999998
$ws.Range("N132").Value = -13399.3334
$ws.Range("H136").Value = 1456.1666
$ws.Range("I136").Value = 1411.5
$ws.Range("K136").Value = 4234.5
$ws.Range("M136").Value = -1684.5

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("H107").Value = 571.4231
$ws.Range("I107").Value = 561.125
$ws.Range("J107").Value = 695
$ws.Range("K107").Value = 561.125
$ws.Range("L107").Value = 695
$ws.Range("M107").Value = 1358.875
$ws.Range("N107").Value = -4535
$ws.Range("H134").Value = 1781.7368
$ws.Range("I134").Value = 1429.5625
$ws.Range("J134").Value = 3660
$ws.Range("K134").Value = 4288.6875
$ws.Range("L134").Value = 10980
$ws.Range("M134").Value = -1753.6875
$ws.Range("N134").Value = -16050
$ws.Range("N69").ClearContents()
$ws.Range("N72").ClearContents()

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 1994.1428
$ws.Range("I16").Value = 1994.1428
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1994.1428
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1707.1428
$ws.Range("H31").Value = 2415.3225
$ws.Range("I31").Value = 2128.8096
$ws.Range("J31").Value = 2562.0732
$ws.Range("K31").Value = 2128.8096
$ws.Range("L31").Value = 2562.0732
$ws.Range("M31").Value = -1833.8096
$ws.Range("N31").Value = -3152.0732
$ws.Range("H34").Value = 2415.3225
$ws.Range("I34").Value = 2128.8096
$ws.Range("J34").Value = 2562.0732
$ws.Range("K34").Value = 2128.8096
$ws.Range("L34").Value = 2562.0732
$ws.Range("M34").Value = -1926.8096
$ws.Range("N34").Value = -2966.0732
$ws.Range("H86").Value = 2666.5
$ws.Range("I86").Value = 2641.1538
$ws.Range("J86").Value = 2732.4
$ws.Range("K86").Value = 2641.1538
$ws.Range("L86").Value = 2732.4
$ws.Range("M86").Value = -1518.1538
$ws.Range("N86").Value = -4978.4
$ws.Range("H88").Value = 20199.8
$ws.Range("J88").Value = 20199.8
$ws.Range("L88").Value = 20199.8
$ws.Range("N88").Value = -21011.8
$ws.Range("H89").Value = 2666.5
$ws.Range("I89").Value = 2641.1538
$ws.Range("J89").Value = 2732.4
$ws.Range("K89").Value = 13205.769
$ws.Range("L89").Value = 13662
$ws.Range("M89").Value = -7589.769
$ws.Range("N89").Value = -24894
$ws.Range("H91").Value = 20199.8
$ws.Range("J91").Value = 20199.8
$ws.Range("L91").Value = 20199.8
$ws.Range("N91").Value = -23007.8
$ws.Range("H99").Value = 8026200
$ws.Range("I99").Value = 10668267
$ws.Range("K99").Value = 10668267
$ws.Range("M99").Value = -10666769
$ws.Range("H107").Value = 583.5238000000001
$ws.Range("I107").Value = 192.75
$ws.Range("J107").Value = 824
$ws.Range("K107").Value = 192.75
$ws.Range("L107").Value = 824
$ws.Range("M107").Value = 1727.25
$ws.Range("N107").Value = -4664
$ws.Range("H113").Value = 1994.1428
$ws.Range("I113").Value = 1994.1428
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1994.1428
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 175.8571999999999
$ws.Range("H126").Value = 8026200
$ws.Range("I126").Value = 10668267
$ws.Range("K126").Value = 32004801
$ws.Range("M126").Value = -32002331
$ws.Range("H132").Value = 3425
$ws.Range("I132").Value = 1555.909
$ws.Range("J132").Value = 5294.091
$ws.Range("K132").Value = 4667.727000000001
$ws.Range("L132").Value = 15882.273
$ws.Range("M132").Value = -2137.727000000001
$ws.Range("N132").Value = -20942.273
$ws.Range("N16").ClearContents()
$ws.Range("N113").ClearContents()

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 4418.9653
$ws.Range("I131").Value = 320
$ws.Range("J131").Value = 5980.476
$ws.Range("K131").Value = 960
$ws.Range("L131").Value = 17941.428
$ws.Range("M131").Value = 4080
$ws.Range("N131").Value = -28021.428
$ws.Range("H132").Value = 1084.375
$ws.Range("J132").Value = 1300
$ws.Range("L132").Value = 11700
$ws.Range("N132").Value = -16760

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H113").Value = 1327.75
$ws.Range("I113").Value = 1327.75
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1327.75
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 842.25
$ws.Range("H132").Value = 3133.8147
$ws.Range("I132").Value = 3298.5293
$ws.Range("K132").Value = 9895.5879
$ws.Range("M132").Value = -7365.5879
$ws.Range("H136").Value = 17766.666
$ws.Range("J136").Value = 17766.666
$ws.Range("L136").Value = 53299.99800000001
$ws.Range("N136").Value = -58399.99800000001
$ws.Range("N113").ClearContents()

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H82").Value = 1746.8
$ws.Range("I82").Value = 923.5714
$ws.Range("J82").Value = 3667.6667
$ws.Range("K82").Value = 923.5714
$ws.Range("L82").Value = 3667.6667
$ws.Range("M82").Value = -562.5714
$ws.Range("N82").Value = -4389.6667
$ws.Range("H85").Value = 1746.8
$ws.Range("I85").Value = 923.5714
$ws.Range("J85").Value = 3667.6667
$ws.Range("K85").Value = 923.5714
$ws.Range("L85").Value = 3667.6667
$ws.Range("M85").Value = 324.4286
$ws.Range("N85").Value = -6163.6667
$ws.Range("H93").Value = 10908.637
$ws.Range("I93").Value = 15898.143
$ws.Range("J93").Value = 2177
$ws.Range("K93").Value = 15898.143
$ws.Range("L93").Value = 2177
$ws.Range("M93").Value = -14650.143
$ws.Range("N93").Value = -4673

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H136").Value = 10757821
$ws.Range("I136").Value = 15152384
$ws.Range("J136").Value = 15556.111
$ws.Range("K136").Value = 45457152
$ws.Range("L136").Value = 46668.333
$ws.Range("M136").Value = -45454602
$ws.Range("N136").Value = -51768.333
